$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bulk Upload")

$ws.Range("D2:D31").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "Proctor08585"
$ws.Cells.Item(2, 2).Value = "Automation08585"
$ws.Cells.Item(2, 3).Value = "proctorautomation08585@gmail.com"
$ws.Cells.Item(2, 4).Value = "08585"

$ws.Cells.Item(3, 1).Value = "Proctor59350"
$ws.Cells.Item(3, 2).Value = "Automation59350"
$ws.Cells.Item(3, 3).Value = "proctorautomation59350@gmail.com"
$ws.Cells.Item(3, 4).Value = "59350"

$ws.Cells.Item(4, 1).Value = "Proctor00477"
$ws.Cells.Item(4, 2).Value = "Automation00477"
$ws.Cells.Item(4, 3).Value = "proctorautomation00477@gmail.com"
$ws.Cells.Item(4, 4).Value = "00477"

$ws.Cells.Item(5, 1).Value = "Proctor21499"
$ws.Cells.Item(5, 2).Value = "Automation21499"
$ws.Cells.Item(5, 3).Value = "proctorautomation21499@gmail.com"
$ws.Cells.Item(5, 4).Value = "21499"

$ws.Cells.Item(6, 1).Value = "Proctor29959"
$ws.Cells.Item(6, 2).Value = "Automation29959"
$ws.Cells.Item(6, 3).Value = "proctorautomation29959@gmail.com"
$ws.Cells.Item(6, 4).Value = "29959"

$ws.Cells.Item(7, 1).Value = "Proctor29280"
$ws.Cells.Item(7, 2).Value = "Automation29280"
$ws.Cells.Item(7, 3).Value = "proctorautomation29280@gmail.com"
$ws.Cells.Item(7, 4).Value = "29280"

$ws.Cells.Item(8, 1).Value = "Proctor22528"
$ws.Cells.Item(8, 2).Value = "Automation22528"
$ws.Cells.Item(8, 3).Value = "proctorautomation22528@gmail.com"
$ws.Cells.Item(8, 4).Value = "22528"

$ws.Cells.Item(9, 1).Value = "Proctor80625"
$ws.Cells.Item(9, 2).Value = "Automation80625"
$ws.Cells.Item(9, 3).Value = "proctorautomation80625@gmail.com"
$ws.Cells.Item(9, 4).Value = "80625"

$ws.Cells.Item(10, 1).Value = "Proctor12961"
$ws.Cells.Item(10, 2).Value = "Automation12961"
$ws.Cells.Item(10, 3).Value = "proctorautomation12961@gmail.com"
$ws.Cells.Item(10, 4).Value = "12961"

$ws.Cells.Item(11, 1).Value = "Proctor36758"
$ws.Cells.Item(11, 2).Value = "Automation36758"
$ws.Cells.Item(11, 3).Value = "proctorautomation36758@gmail.com"
$ws.Cells.Item(11, 4).Value = "36758"

$ws.Cells.Item(12, 1).Value = "Proctor50433"
$ws.Cells.Item(12, 2).Value = "Automation50433"
$ws.Cells.Item(12, 3).Value = "proctorautomation50433@gmail.com"
$ws.Cells.Item(12, 4).Value = "50433"

$ws.Cells.Item(13, 1).Value = "Proctor48286"
$ws.Cells.Item(13, 2).Value = "Automation48286"
$ws.Cells.Item(13, 3).Value = "proctorautomation48286@gmail.com"
$ws.Cells.Item(13, 4).Value = "48286"

$ws.Cells.Item(14, 1).Value = "Proctor45318"
$ws.Cells.Item(14, 2).Value = "Automation45318"
$ws.Cells.Item(14, 3).Value = "proctorautomation45318@gmail.com"
$ws.Cells.Item(14, 4).Value = "45318"

$ws.Cells.Item(15, 1).Value = "Proctor79288"
$ws.Cells.Item(15, 2).Value = "Automation79288"
$ws.Cells.Item(15, 3).Value = "proctorautomation79288@gmail.com"
$ws.Cells.Item(15, 4).Value = "79288"

$ws.Cells.Item(16, 1).Value = "Proctor39503"
$ws.Cells.Item(16, 2).Value = "Automation39503"
$ws.Cells.Item(16, 3).Value = "proctorautomation39503@gmail.com"
$ws.Cells.Item(16, 4).Value = "39503"

$ws.Cells.Item(17, 1).Value = "Proctor48235"
$ws.Cells.Item(17, 2).Value = "Automation48235"
$ws.Cells.Item(17, 3).Value = "proctorautomation48235@gmail.com"
$ws.Cells.Item(17, 4).Value = "48235"

$ws.Cells.Item(18, 1).Value = "Proctor15585"
$ws.Cells.Item(18, 2).Value = "Automation15585"
$ws.Cells.Item(18, 3).Value = "proctorautomation15585@gmail.com"
$ws.Cells.Item(18, 4).Value = "15585"

$ws.Cells.Item(19, 1).Value = "Proctor10556"
$ws.Cells.Item(19, 2).Value = "Automation10556"
$ws.Cells.Item(19, 3).Value = "proctorautomation10556@gmail.com"
$ws.Cells.Item(19, 4).Value = "10556"

$ws.Cells.Item(20, 1).Value = "Proctor29091"
$ws.Cells.Item(20, 2).Value = "Automation29091"
$ws.Cells.Item(20, 3).Value = "proctorautomation29091@gmail.com"
$ws.Cells.Item(20, 4).Value = "29091"

$ws.Cells.Item(21, 1).Value = "Proctor73299"
$ws.Cells.Item(21, 2).Value = "Automation73299"
$ws.Cells.Item(21, 3).Value = "proctorautomation73299@gmail.com"
$ws.Cells.Item(21, 4).Value = "73299"

$ws.Cells.Item(22, 1).Value = "Proctor94128"
$ws.Cells.Item(22, 2).Value = "Automation94128"
$ws.Cells.Item(22, 3).Value = "proctorautomation94128@gmail.com"
$ws.Cells.Item(22, 4).Value = "94128"

$ws.Cells.Item(23, 1).Value = "Proctor88742"
$ws.Cells.Item(23, 2).Value = "Automation88742"
$ws.Cells.Item(23, 3).Value = "proctorautomation88742@gmail.com"
$ws.Cells.Item(23, 4).Value = "88742"

$ws.Cells.Item(24, 1).Value = "Proctor41863"
$ws.Cells.Item(24, 2).Value = "Automation41863"
$ws.Cells.Item(24, 3).Value = "proctorautomation41863@gmail.com"
$ws.Cells.Item(24, 4).Value = "41863"

$ws.Cells.Item(25, 1).Value = "Proctor29685"
$ws.Cells.Item(25, 2).Value = "Automation29685"
$ws.Cells.Item(25, 3).Value = "proctorautomation29685@gmail.com"
$ws.Cells.Item(25, 4).Value = "29685"

$ws.Cells.Item(26, 1).Value = "Proctor11339"
$ws.Cells.Item(26, 2).Value = "Automation11339"
$ws.Cells.Item(26, 3).Value = "proctorautomation11339@gmail.com"
$ws.Cells.Item(26, 4).Value = "11339"

$ws.Cells.Item(27, 1).Value = "Proctor27409"
$ws.Cells.Item(27, 2).Value = "Automation27409"
$ws.Cells.Item(27, 3).Value = "proctorautomation27409@gmail.com"
$ws.Cells.Item(27, 4).Value = "27409"

$ws.Cells.Item(28, 1).Value = "Proctor07476"
$ws.Cells.Item(28, 2).Value = "Automation07476"
$ws.Cells.Item(28, 3).Value = "proctorautomation07476@gmail.com"
$ws.Cells.Item(28, 4).Value = "07476"

$ws.Cells.Item(29, 1).Value = "Proctor61224"
$ws.Cells.Item(29, 2).Value = "Automation61224"
$ws.Cells.Item(29, 3).Value = "proctorautomation61224@gmail.com"
$ws.Cells.Item(29, 4).Value = "61224"

$ws.Cells.Item(30, 1).Value = "Proctor32501"
$ws.Cells.Item(30, 2).Value = "Automation32501"
$ws.Cells.Item(30, 3).Value = "proctorautomation32501@gmail.com"
$ws.Cells.Item(30, 4).Value = "32501"

$ws.Cells.Item(31, 1).Value = "Proctor85926"
$ws.Cells.Item(31, 2).Value = "Automation85926"
$ws.Cells.Item(31, 3).Value = "proctorautomation85926@gmail.com"
$ws.Cells.Item(31, 4).Value = "85926"

$ws.Range("D2:D31").Style = "Normal"